# "Generate Report for Handoff"
# Adds two new localized source files
#   4c91c73a-5999-4d30-87c7-516f8f0137fa.md  (status: Ready for handoff)
#   c738fd0b-fa15-4575-ad70-ec3434746afe.md  (status: Ready for handoff)
# to the Overview/zh-cn/de-de report sheets, pushing the existing
# ".localization-config" row down below them.

$wb = $excel.ActiveWorkbook

$mdRepoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/6575079830b878e0aecfa456478e4dfb7b10d350"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eedee25d9c62ff00ad66720b3fb9aaa8a3c91f21/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68e045f72d494d1714cee9bae17beece81e5ca82/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() on any range clears every hyperlink on the sheet in
# this engine, so start clean and re-add all five (including the three that
# are unchanged) in the final left-to-right / top-to-bottom order.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Value2 = "In Translation"
$ws.Range("C2").Value2 = "In Translation"
$ws.Range("B3").Value2 = "In Translation"
$ws.Range("C3").Value2 = "In Translation"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "Ready for handoff"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "Ready for handoff"
$ws.Range("B6").Value2 = "Not to be localized"
$ws.Range("C6").Value2 = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdRepoBase/e2e/93a70da4-373d-46a8-9c30-01a8ea25c1e9.md", "", "", "93a70da4-373d-46a8-9c30-01a8ea25c1e9.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdRepoBase/e2e/f6501bc2-a484-4f7c-98bf-aaf53afd0c17.md", "", "", "f6501bc2-a484-4f7c-98bf-aaf53afd0c17.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdRepoBase/e2e/4c91c73a-5999-4d30-87c7-516f8f0137fa.md", "", "", "4c91c73a-5999-4d30-87c7-516f8f0137fa.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdRepoBase/e2e/c738fd0b-fa15-4575-ad70-ec3434746afe.md", "", "", "c738fd0b-fa15-4575-ad70-ec3434746afe.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdRepoBase/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Value2 = "In Translation"
$ws.Range("D2").Value2 = "2016-01-27 07:51:29"
$ws.Range("G2").Value2 = "0001-01-01 00:00:00"
$ws.Range("H2").Value2 = "Include"

$ws.Range("B3").Value2 = "In Translation"
$ws.Range("D3").Value2 = "2016-01-27 07:51:29"
$ws.Range("G3").Value2 = "0001-01-01 00:00:00"
$ws.Range("H3").Value2 = "Include"

$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("D4").Value2 = "2016-01-27 07:53:05"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Include"

$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("D5").Value2 = "2016-01-27 07:53:05"
$ws.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws.Range("H5").Value2 = "Include"

$ws.Range("B6").Value2 = "Not to be localized"
$ws.Range("D6").Value2 = "0001-01-01 00:00:00"
$ws.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws.Range("H6").Value2 = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdRepoBase/e2e/93a70da4-373d-46a8-9c30-01a8ea25c1e9.md", "", "", "93a70da4-373d-46a8-9c30-01a8ea25c1e9.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "$zhHandoffBase/93a70da4-373d-46a8-9c30-01a8ea25c1e9.074e41af346f0c6b213c8c6715aedf7e5d069ede.zh-cn.xlf", "", "", "93a70da4-373d-46a8-9c30-01a8ea25c1e9.074e41af346f0c6b213c8c6715aedf7e5d069ede.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdRepoBase/e2e/f6501bc2-a484-4f7c-98bf-aaf53afd0c17.md", "", "", "f6501bc2-a484-4f7c-98bf-aaf53afd0c17.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "$zhHandoffBase/f6501bc2-a484-4f7c-98bf-aaf53afd0c17.898a68cbc526cbe8b8d246470ece471a4fbf2c6b.zh-cn.xlf", "", "", "f6501bc2-a484-4f7c-98bf-aaf53afd0c17.898a68cbc526cbe8b8d246470ece471a4fbf2c6b.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdRepoBase/e2e/4c91c73a-5999-4d30-87c7-516f8f0137fa.md", "", "", "4c91c73a-5999-4d30-87c7-516f8f0137fa.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "$zhHandoffBase/4c91c73a-5999-4d30-87c7-516f8f0137fa.fd98bb5f7345ddbc51c7965d250a58113e9b1a6b.zh-cn.xlf", "", "", "4c91c73a-5999-4d30-87c7-516f8f0137fa.fd98bb5f7345ddbc51c7965d250a58113e9b1a6b.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdRepoBase/e2e/c738fd0b-fa15-4575-ad70-ec3434746afe.md", "", "", "c738fd0b-fa15-4575-ad70-ec3434746afe.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "$zhHandoffBase/c738fd0b-fa15-4575-ad70-ec3434746afe.b40dc611742437118fcd6611e8414ef3cb662c6f.zh-cn.xlf", "", "", "c738fd0b-fa15-4575-ad70-ec3434746afe.b40dc611742437118fcd6611e8414ef3cb662c6f.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdRepoBase/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Value2 = "In Translation"
$ws.Range("D2").Value2 = "2016-01-27 07:51:54"
$ws.Range("G2").Value2 = "0001-01-01 00:00:00"
$ws.Range("H2").Value2 = "Include"

$ws.Range("B3").Value2 = "In Translation"
$ws.Range("D3").Value2 = "2016-01-27 07:51:54"
$ws.Range("G3").Value2 = "0001-01-01 00:00:00"
$ws.Range("H3").Value2 = "Include"

$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("D4").Value2 = "2016-01-27 07:53:16"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Include"

$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("D5").Value2 = "2016-01-27 07:53:16"
$ws.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws.Range("H5").Value2 = "Include"

$ws.Range("B6").Value2 = "Not to be localized"
$ws.Range("D6").Value2 = "0001-01-01 00:00:00"
$ws.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws.Range("H6").Value2 = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdRepoBase/e2e/93a70da4-373d-46a8-9c30-01a8ea25c1e9.md", "", "", "93a70da4-373d-46a8-9c30-01a8ea25c1e9.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "$deHandoffBase/93a70da4-373d-46a8-9c30-01a8ea25c1e9.074e41af346f0c6b213c8c6715aedf7e5d069ede.de-de.xlf", "", "", "93a70da4-373d-46a8-9c30-01a8ea25c1e9.074e41af346f0c6b213c8c6715aedf7e5d069ede.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdRepoBase/e2e/f6501bc2-a484-4f7c-98bf-aaf53afd0c17.md", "", "", "f6501bc2-a484-4f7c-98bf-aaf53afd0c17.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "$deHandoffBase/f6501bc2-a484-4f7c-98bf-aaf53afd0c17.898a68cbc526cbe8b8d246470ece471a4fbf2c6b.de-de.xlf", "", "", "f6501bc2-a484-4f7c-98bf-aaf53afd0c17.898a68cbc526cbe8b8d246470ece471a4fbf2c6b.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdRepoBase/e2e/4c91c73a-5999-4d30-87c7-516f8f0137fa.md", "", "", "4c91c73a-5999-4d30-87c7-516f8f0137fa.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "$deHandoffBase/4c91c73a-5999-4d30-87c7-516f8f0137fa.fd98bb5f7345ddbc51c7965d250a58113e9b1a6b.de-de.xlf", "", "", "4c91c73a-5999-4d30-87c7-516f8f0137fa.fd98bb5f7345ddbc51c7965d250a58113e9b1a6b.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdRepoBase/e2e/c738fd0b-fa15-4575-ad70-ec3434746afe.md", "", "", "c738fd0b-fa15-4575-ad70-ec3434746afe.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "$deHandoffBase/c738fd0b-fa15-4575-ad70-ec3434746afe.b40dc611742437118fcd6611e8414ef3cb662c6f.de-de.xlf", "", "", "c738fd0b-fa15-4575-ad70-ec3434746afe.b40dc611742437118fcd6611e8414ef3cb662c6f.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdRepoBase/.localization-config", "", "", ".localization-config")

Write-Output "Report regenerated for handoff."
